$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range("G3").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G4").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G5").Value = 'servinaz@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hananragab@med.asu.edu.eg'
$ws.Range("G6").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G7").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G8").Value = 'asmaa.reda@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, backup@backdoor.com, majorelle.magdy@med.asu.edu.eg'
$ws.Range("G9").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg'
$ws.Range("G10").Value = 'heba@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, gehanadel@med.asu.edu.eg, drsarahmw@gmail.com'
$ws.Range("G11").Value = 'Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range("G13").Value = 'wessam.atef@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G17").Value = 'yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, dina.adel@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg'
$ws.Range("G24").Value = 'yasmin.m.senosy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg'
$ws.Range("G25").Value = 'abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg'
$ws.Range("G26").Value = 'gigi@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg'
$ws.Range("G27").Value = 'yasmin.m.senosy@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range("G28").Value = 'abdullah.elagrody@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range("G29").Value = 'ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, eman.samir@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G30").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G31").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G32").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G33").Value = 'servinaz@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hananragab@med.asu.edu.eg'
$ws.Range("G34").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G35").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G36").Value = 'asmaa.reda@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, backup@backdoor.com, majorelle.magdy@med.asu.edu.eg'
$ws.Range("G37").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg'
$ws.Range("G38").Value = 'heba@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, gehanadel@med.asu.edu.eg, drsarahmw@gmail.com'
$ws.Range("G39").Value = 'Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range("G41").Value = 'wessam.atef@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G45").Value = 'yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, dina.adel@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg'
$ws.Range("G52").Value = 'yasmin.m.senosy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg'
$ws.Range("G53").Value = 'abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg'
$ws.Range("G54").Value = 'gigi@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg'
$ws.Range("G55").Value = 'yasmin.m.senosy@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range("G56").Value = 'abdullah.elagrody@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range("G57").Value = 'ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, eman.samir@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
